# Updates the cryptos list values (Price / Volume(1h)) to reflect a new
# data snapshot, and swaps the EnergySwap / Decentraland rows (46 <-> 47).
#
# Price (column D) and Volume(1h) (column E) are stored as plain text in
# this sheet (e.g. '1.001', '28.575.57'), so each cell is forced to Text
# number format before the value is written (and reset back to the default
# 'Normal' style afterwards) to stop Excel from auto-converting numeric-
# looking strings into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '28.575.57'
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.34%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '1.848.65'
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.93%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.01%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '333.64'
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = '  +2.23%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.99%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.4658'
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.44%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = '0.3919'
$c.Style = "Normal"
$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.79%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '46.51'
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.50%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.07914'
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = '  -4.02%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '0.9850'
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.65%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '22.21'
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = '  -5.99%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '1.994.63'
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = '  +4.15%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '5.847'
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.95%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '7.024'
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.58%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '0.06885'
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value = '  +0.08%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '87.81'
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = '  -4.36%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.03%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '0.00001006'
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.23%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '17.09'
$c.Style = "Normal"
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.00%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = '  -0.81%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '28.608.14'
$c.Style = "Normal"
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.26%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '5.402'
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = '  -5.01%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = '  -5.24%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '2.213.63'
$c.Style = "Normal"
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = '  +4.10%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '2.128'
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.66%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '153.22'
$c.Style = "Normal"
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = '  -1.83%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '19.41'
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.18%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '6.119'
$c.Style = "Normal"
$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = '  -5.63%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '2.012'
$c.Style = "Normal"
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = '  -4.18%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '117.63'
$c.Style = "Normal"
$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.56%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '0.9830'
$c.Style = "Normal"
$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.33%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '0.09428'
$c.Style = "Normal"
$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.26%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '5.367'
$c.Style = "Normal"
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = '  -4.78%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '3.484'
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.02%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '1.347'
$c.Style = "Normal"
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.16%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '0.06150'
$c.Style = "Normal"
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.53%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '0.02203'
$c.Style = "Normal"
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = '  -4.14%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '1.159'
$c.Style = "Normal"
$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.79%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = '  -4.28%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '7.611'
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.69%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '10.13'
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = '  -6.23%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '0.1797'
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.93%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '2.368'
$c.Style = "Normal"
$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = '  -4.22%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '1.250'
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.43%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '0.07156'
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = '  -4.67%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '1.907'
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = '  -2.40%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '114.11'
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = '  -4.18%  '
$c.Style = "Normal"

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '42.75'
$c.Style = "Normal"
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = '  +1.74%  '
$c.Style = "Normal"

# Rows 46/47: EnergySwap and Decentraland swap order, with new Price/Volume values.
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '11.88'
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = '  -4.58%  '
$c.Style = "Normal"

$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '0.5397'
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = '  -3.31%  '
$c.Style = "Normal"

